$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.5908239700374532
$ws1.Range("C2").Value = 0.552319309600863
$ws1.Range("D2").Value = 0.9588014981273408
$ws1.Range("E2").Value = 0.7008898015058179
$ws1.Range("F2").Value = 0.8357819131570355
$ws1.Range("G2").Value = 0.9324087693493031
$ws1.Range("H2").Value = 0.7935656272356184
$ws1.Range("I2").Value = 512
$ws1.Range("J2").Value = 415
$ws1.Range("K2").Value = 119
$ws1.Range("L2").Value = 22

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# row 2 (label "0")
$ws2.Range("B2").Value = 0.8439716312056738
$ws2.Range("C2").Value = 0.2228464419475655
$ws2.Range("D2").Value = 0.3525925925925926

# row 3 (label "1")
$ws2.Range("B3").Value = 0.552319309600863
$ws2.Range("C3").Value = 0.9588014981273408
$ws2.Range("D3").Value = 0.7008898015058179

# row 4 (label "accuracy")
$ws2.Range("B4").Value = 0.5908239700374532
$ws2.Range("C4").Value = 0.5908239700374532
$ws2.Range("D4").Value = 0.5908239700374532
$ws2.Range("E4").Value = 0.5908239700374532

# row 5 (label "macro avg")
$ws2.Range("B5").Value = 0.6981454704032684
$ws2.Range("C5").Value = 0.5908239700374531
$ws2.Range("D5").Value = 0.5267411970492053

# row 6 (label "weighted avg")
$ws2.Range("B6").Value = 0.6981454704032684
$ws2.Range("C6").Value = 0.5908239700374532
$ws2.Range("D6").Value = 0.5267411970492053

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# row 2 (Actual 0)
$ws3.Range("B2").Value = 119
$ws3.Range("C2").Value = 415

# row 3 (Actual 1)
$ws3.Range("B3").Value = 22
$ws3.Range("C3").Value = 512
